# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.986.57"
$ws.Range("D3").Value = "1.859.16"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "312.06"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3830"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").Value = "0.08218"
$ws.Range("E9").Value = "  -10.19%  "
$ws.Range("D10").Value = "1.109"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").Value = "41.52"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "6.189"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").Value = "20.55"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "1.858.78"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "7.241"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.60"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "0.06645"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").Value = "17.69"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "6.004"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("D23").Value = "28.010.43"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").Value = "2.244"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").Value = "2.072.23"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").Value = "2.509"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").Value = "157.76"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "20.44"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").Value = "124.57"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.030"
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("D33").Value = "5.941"
$ws.Range("E33").Value = "  +5.86%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "9.361"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("D36").Value = "0.02414"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06490"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").Value = "0.2173"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "0.6531"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").Value = "1.194"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").Value = "4.995"
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("D42").Value = "1.217"
$ws.Range("E42").Value = "  -3.36%  "
$ws.Range("D43").Value = "11.17"
$ws.Range("E43").Value = "  -3.78%  "
$ws.Range("D44").Value = "0.6148"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("D45").Value = "12.97"
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "3.657"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "2.004"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("D50").Value = "120.47"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "78.02"
$ws.Range("E51").Value = "  -2.20%  "
